{"js": "// Update the business address printed in the page header, and drop the\n// stale \"_GoBack\" bookmark that Word leaves behind at the last edit\n// position (both are what the commit \"Updated address on the docs\" did).\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst oldAddress = \"1175 Whitehawk Trail, Lawrenceville, GA 30043\";\nconst newAddress = \"1969 Mapmaker Drive, Dacula, GA 30019\";\n\n// The address lives in the primary (default) header of the first section.\nconst header = sections.items[0].getHeader(\"Primary\");\nconst found = header.search(oldAddress, { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\n\nif (found.items.length > 0) {\n  found.items[0].insertText(newAddress, \"Replace\");\n}\n\n// Remove the leftover \"_GoBack\" bookmark from the body.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Update the business address printed in the page header, and drop the\n# stale \"_GoBack\" bookmark that Word leaves behind at the last edit\n# position (both are what the commit \"Updated address on the docs\" did).\n\n$d = $word.ActiveDocument\n\n$oldAddress = \"1175 Whitehawk Trail, Lawrenceville, GA 30043\"\n$newAddress = \"1969 Mapmaker Drive, Dacula, GA 30019\"\n\n# The address lives in the primary (default) header of the first section.\n$header = $d.Sections.Item(1).Headers.Item(1)\n$find = $header.Range.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldAddress\n$find.Replacement.Text = $newAddress\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Remove the leftover \"_GoBack\" bookmark from the body.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n"}
